$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers for the additional measured sessions (tabs) columns
$ws.Range("D1").Value = "freq_2"
$ws.Range("E1").Value = "help_2"

# Update the selection to match the new active cell
$ws.Range("D2").Select()
